$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 / Row 10: the two records effectively swap their
# prolificid + name (shared-string reorder in the source diff),
# and also swap their "race" (C) and "gender" (F) values; the
# realeffort (G) scores are refreshed for every data row below. ---

$ws.Range("D9").Value  = "60bd88b8fc436774352f53b9"
$ws.Range("E9").Value  = "Annes"
$ws.Range("C9").Value  = 3
$ws.Range("F9").Value  = "female"

$ws.Range("D10").Value = "5c27de12a2b00a00018b2c16"
$ws.Range("E10").Value = "Ankai"
$ws.Range("C10").Value = 0
$ws.Range("F10").Value = "male"

# --- Refresh the "realeffort" (G column) scores for all 24 data rows ---

$ws.Range("G2").Value  = 11.12356095231806
$ws.Range("G3").Value  = 10.07657103797102
$ws.Range("G4").Value  = 8.469824362969149
$ws.Range("G5").Value  = 8.252431536799262
$ws.Range("G6").Value  = 7.094657342882389
$ws.Range("G7").Value  = 6.351992923050718
$ws.Range("G8").Value  = 6.111943368614604
$ws.Range("G9").Value  = 5.465857846036377
$ws.Range("G10").Value = 5.411049145544538
$ws.Range("G11").Value = 4.011565163053068
$ws.Range("G12").Value = 2.037811163075423
$ws.Range("G13").Value = 0.2234880371687213
$ws.Range("G14").Value = 8.487299773058886
$ws.Range("G15").Value = 8.17703893189862
$ws.Range("G16").Value = 7.389685286561035
$ws.Range("G17").Value = 7.23768709675885
$ws.Range("G18").Value = 6.343048961041616
$ws.Range("G19").Value = 6.040248810237871
$ws.Range("G20").Value = 5.132481009527025
$ws.Range("G21").Value = 3.39083832478149
$ws.Range("G22").Value = 3.318523148442117
$ws.Range("G23").Value = 2.176093329063497
$ws.Range("G24").Value = 1.474100674863718
$ws.Range("G25").Value = 0.05652031328557322
